$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.620.47"
$ws.Range("E2").Value = "  +4.17%  "
$ws.Range("D3").Value = "1.746.13"
$ws.Range("E3").Value = "  +4.50%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "247.32"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.42%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.4802"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.47%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.2698"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.78%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.06263"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "1.743.99"
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07112"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("E12").Value = "  +6.29%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.6171"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.59%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "4.511"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.09%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "77.25"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D17").Value = "26.620.95"
$ws.Range("E17").Value = "  +4.20%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.10%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.000006900"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.02%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "11.72"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").Value = "1.970.56"
$ws.Range("E21").Value = "  +4.59%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "4.645"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +4.53%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "8.862"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.25%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "5.346"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.29%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "136.15"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.36%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "15.46"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.89%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "1.823"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.83%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "1.413"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "107.58"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.70%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "4.024"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.59%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "3.776"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +3.43%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.07899"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.82%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.04580"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +8.08%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "2.614"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.29%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.57%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.6364"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.63%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.9487"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +9.62%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "114.15"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +18.33%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "2.465"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -5.06%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "1.984"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.21%  "
$ws.Range("E41").Value = "  +0.40%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.01512"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.24%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "5.663"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +16.00%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.3915"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.93%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "6.725"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +8.21%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.1203"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.74%  "
$ws.Range("E47").Value = "  +1.24%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "7.945"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +6.99%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "30.90"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("E50").Value = "  +4.29%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.3457"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.56%  "
